$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-7 from 45175 to 45183
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
